$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "groundwater recharge" row (row 12) entirely. This shifts the
# "nanofiltration" row (old row 13) up to row 12, and since the
# "groundwater recharge" string becomes unused it drops out of the shared
# strings table, shifting later string indices down by one.
$ws.Rows(12).Delete()

# Re-assign the formula that got shifted up so it is re-serialized in its
# normal form (avoids a spurious _xlfn.SINGLE implicit-intersection wrapper
# that recalculation-after-row-shift can introduce).
$ws.Range("I12").Formula = $ws.Range("I12").Formula

# Match the resulting selection/cursor position left behind by the edit.
$ws.Range("C15").Select() | Out-Null
